$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 4.8
$ws.Range("F3").Value = 1.44
$ws.Range("K3").Value = 4.8
$ws.Range("F4").Value = 1.66
$ws.Range("G4").Value = 1.67
$ws.Range("Q4").Value = 1.81
$ws.Range("AN4").Value = 8.4
